$wb = $excel.ActiveWorkbook
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $after)
$ws4.Name = "Access_Permission"

$ws4.Range("A1").Value = "Groups"
$ws4.Range("A2").Value = "GroupName"
$ws4.Range("B2").Value = "Group Description"
$ws4.Range("A3").Value = "Groupssss"
$ws4.Range("B3").Value = "This is Group Description"

$ws4.Range("A5").Value = "Roles"
$ws4.Range("A6").Value = "Role Name"
$ws4.Range("B6").Value = "Select Group"
$ws4.Range("C6").Value = "Role Description"
$ws4.Range("A7").Value = "Rolessss"
$ws4.Range("B7").Value = "Super Admin"
$ws4.Range("C7").Value = "This is Role Description"

$ws4.Range("A9").Value = "Users"
$ws4.Range("A10").Value = "LoginId"
$ws4.Range("B10").Value = "FirstName"
$ws4.Range("C10").Value = "LastName"
$ws4.Range("D10").Value = "EmailId"
$ws4.Range("E10").Value = "Department"
$ws4.Range("F10").Value = "PhoneNumber"

$ws4.Range("A11").Value = "User2810"
$ws4.Range("B11").Value = "User Fname"
$ws4.Range("C11").Value = "User Lname"
$ws4.Range("D11").Value = "test@gmail.com"
$ws4.Range("E11").Value = "Testing"
$ws4.Range("F11").NumberFormat = "@"
$ws4.Range("F11").Value = "8881212888"

Write-Host "values done"
